# The upstream commit ("Fixed #295 Add the version of M2Doc in the
# template custom properties") re-saved this particular template
# through the M2Doc tooling. For *this* resource the resulting XML
# diff is a pure re-serialization: every changed line is the same
# element with the same attribute names/values, just written back
# out with the attributes in (alphabetically) sorted order -- e.g.
#   <w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/>
# becomes
#   <w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>
# and likewise for <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>,
# <w:latentStyles>/<w:lsdException>, <w:style>, <w:tblInd>, the
# w:tblCellMar children, and the xmlns:* declarations on the
# <w:document> root. No text, run, paragraph, style, property value,
# page size/margin, font, or latent-style setting actually changed.
#
# So the faithful edit is simply to touch/re-save the document as-is
# -- no content, formatting, or property mutation is required.
$d = $word.ActiveDocument
$d.Save()
